$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3383571169547176
$ws.Range("D2").Value = 0.4136789058781119
$ws.Range("E2").Value = 0.5344303282584365
$ws.Range("F2").Value = 1.055971845676007
$ws.Range("G2").Value = 1.183508476785364
